$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Copy formatting from row 2's already-filled cells so rows 3-5 match style
$ws.Range("G2").Copy()
$ws.Range("G3:G5").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("L2").Copy()
$ws.Range("L3:L5").PasteSpecial(-4122)  # xlPasteFormats

# Fill in the Young's modulus (E) values for rows 3-5 in column G
$ws.Range("G3").Value = 300
$ws.Range("G4").Value = 300
$ws.Range("G5").Value = 300

# Fill in the E_type values for rows 3-5 in column L
$ws.Range("L3").Value = 9
$ws.Range("L4").Value = 9
$ws.Range("L5").Value = 9

# Update the active selection to G4 as in the edited workbook
$ws.Range("G4").Select()
